# "Generate Report for Handback" - refresh the handoff/handback timestamps
# recorded for the zh-cn and de-de handback status reports.

$wb = $excel.ActiveWorkbook

# zh-cn sheet: Correspond Handoff Datetime (E2) / Correspond Handback DateTime (H2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-23 23:14:12"
$wsZhCn.Range("H2").Value = "2016-03-23 23:14:41"

# de-de sheet: Correspond Handoff Datetime (E2) / Correspond Handback DateTime (H2)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-23 23:14:21"
$wsDeDe.Range("H2").Value = "2016-03-23 23:14:48"
